$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$json = @'
{
"公司名称": "云账房",
"地址": "地址：南站绿地之窗",
"部门": "税务事业部",
"岗位": "前端",
"上班时间": "9:00-18:00",
"午休时长": "1.5h",
"加班情况": "发版的时候，会加到10点。上面领导抓得紧的话，会要求每周两天加到8点，没事也要干坐着",
"公积金比例": "总薪资 * 0.8 * 0.6 * 10%",
"年终奖": "去年只发了半个月",
"试用期工资": "三个月，薪资不打折(可跟HR谈)",
"工位电脑情况": "网吧工位，一个台式主机，两个24寸1080P显示器。椅子最烂了，坐着咯吱咯吱响",
"年假": "五天年假，按入职日期折算，每两个月发放一天",
"打卡情况": "钉钉打卡，每个月五次补卡机会",
"其他备注": "调薪需要答辩，还不一定有名额，通过率也不是很高，过了涨幅在1-2千左右。还不是每年都有答辩，看公司情况"
}
'@

$ws.Range("P7").Value = $json
$ws.Range("Q7").Value = "2022-02-11 01:20:28"

# Keep row height consistent with an unmodified row (avoid Excel's
# automatic row-height growth from the newly-added multi-line text).
$ws.Rows.Item(7).AutoFit()
